# Updated symbol list on Sat Jan 28 21:36:30 UTC 2023 with GitHub Actions
# Refresh crypto price/volume figures in the "cryptos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched Price/Volume cells to Text so the new values are stored
# as literal strings (matching the existing inline-string cells) instead of
# being auto-parsed into numbers/percentages by Excel.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "305.75"
$ws.Range("E2").Value = "-0.65%"
$ws.Range("D3").Value = "38.90"
$ws.Range("E3").Value = "7.38%"
$ws.Range("D4").Value = "5.105"
$ws.Range("E4").Value = "0.91%"
$ws.Range("D5").Value = "0.08079"
$ws.Range("E5").Value = "-0.44%"
$ws.Range("D6").Value = "1.927"
$ws.Range("E6").Value = "-4.02%"
$ws.Range("D7").Value = "4.202"
$ws.Range("E7").Value = "0.91%"
$ws.Range("D8").Value = "8.036"
$ws.Range("E8").Value = "2.27%"
$ws.Range("D9").Value = "0.9279"
$ws.Range("E9").Value = "0.11%"
$ws.Range("D10").Value = "0.1455"
$ws.Range("E10").Value = "-2.27%"
$ws.Range("D11").Value = "0.1915"
$ws.Range("E11").Value = "-0.69%"
$ws.Range("D12").Value = "0.09027"
$ws.Range("E12").Value = "-0.56%"
$ws.Range("D13").Value = "0.03512"
$ws.Range("E13").Value = "-0.05%"
$ws.Range("D14").Value = "0.09773"
$ws.Range("E14").Value = "-1.16%"
$ws.Range("D15").Value = "0.001396"
$ws.Range("E15").Value = "-1.78%"
$ws.Range("D16").Value = "0.005828"
$ws.Range("E16").Value = "-4.48%"
$ws.Range("D17").Value = "3.782"
$ws.Range("E17").Value = "-1.52%"
$ws.Range("E18").Value = "-0.76%"
$ws.Range("E19").Value = "-0.78%"
$ws.Range("E20").Value = "2.60%"
$ws.Range("D21").Value = "4.698"
$ws.Range("E21").Value = "-2.42%"
$ws.Range("D22").Value = "0.2419"
$ws.Range("E22").Value = "3.16%"
$ws.Range("D23").Value = "0.04378"
$ws.Range("E23").Value = "-0.19%"
$ws.Range("D24").Value = "0.001238"
$ws.Range("E24").Value = "0.51%"
$ws.Range("D25").Value = "0.004271"
$ws.Range("E25").Value = "2.10%"
$ws.Range("E26").Value = "0.11%"
$ws.Range("D39").Value = "0.02024"
$ws.Range("E39").Value = "-1.20%"
$ws.Range("D40").Value = "0.05051"
$ws.Range("E40").Value = "-1.38%"
$ws.Range("D41").Value = "0.007533"
$ws.Range("E41").Value = "0.70%"
$ws.Range("D42").Value = "0.009798"
$ws.Range("E42").Value = "-2.06%"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").Value = "-0.83%"
$ws.Range("D45").Value = "0.009941"
$ws.Range("E45").Value = "0.78%"
$ws.Range("E46").Value = "-1.81%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.13%"
$ws.Range("D48").Value = "0.002876"
$ws.Range("D49").Value = "0.001805"
$ws.Range("E49").Value = "12.65%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.13%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.13%"

# Restore the cells to the workbook's default (un-styled) appearance now
# that the literal text has been committed.
$valueRange.Style = "Normal"
